# Update tax credit values based on actual discount rates (#116)
#
# production_tax_credits: rows 2-9 get per-technology "actual" discount
# rates (column E) instead of the flat 7% placeholder; the recalculated
# CRF / CRF-real / levelized-tax-credit columns (F, G, I) follow via the
# existing formulas. The whole discount-rate column (E2:E18) is reformatted
# to "0.000" (rows 10-18 keep their original 0.07 / 0.1 values but pick up
# the new number format too).
#
# Selections/active-tab move from production_tax_credits!E2 to
# investment_tax_credits!G9 (investment_tax_credits becomes the active
# sheet), while production_tax_credits itself remembers a selection of I2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # production_tax_credits
$ws2 = $wb.Worksheets.Item(2)   # investment_tax_credits

# --- production_tax_credits: new discount rates for rows 2-9 -----------
$ws1.Range("E2").Value2 = 0.0320752195121951
$ws1.Range("E3").Value2 = 0.0564730561021376
$ws1.Range("E4").Value2 = 0.051519516365778
$ws1.Range("E5").Value2 = 0.0538314857296738
$ws1.Range("E6").Value2 = 0.0438437157985803
$ws1.Range("E7").Value2 = 0.0519007613262936
$ws1.Range("E8").Value2 = 0.0515227657596506
$ws1.Range("E9").Value2 = 0.0515227657596506

# The whole discount-rate column now shows three decimals.
$ws1.Range("E2:E18").NumberFormat = "0.000"

# --- selection bookkeeping ----------------------------------------------
# Leave a lingering selection of I2 on production_tax_credits ...
[void]$ws1.Range("I2").Select()

# ... then make investment_tax_credits the active sheet with G9 selected.
[void]$ws2.Select()
[void]$ws2.Range("G9").Select()
